$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K: 2023 data, extending the table that currently ends at column J ---

# Header (row 3) - "2023"
$c = $ws.Range("K3")
$c.Value = 2023
$c.Borders.Item(8).LineStyle = 1
$c.Borders.Item(10).LineStyle = 1
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.ColorIndex = 1
$c.Interior.ThemeColor = 0
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

# "Number of employees" (row 4)
$c = $ws.Range("K4")
$c.Value = 10753
$c.Borders.Item(8).LineStyle = 1
$c.Borders.Item(10).LineStyle = 1
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.ColorIndex = 1
$c.Interior.ThemeColor = 0
$c.VerticalAlignment = -4108
$c.NumberFormat = "#\ ##0"

# "Of which: Women" (row 5)
$c = $ws.Range("K5")
$c.Value = 4558
$c.Borders.Item(10).LineStyle = 1
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.ColorIndex = 1
$c.Interior.ThemeColor = 0
$c.VerticalAlignment = -4108
$c.NumberFormat = "#\ ##0"

# "Men" (row 6)
$c = $ws.Range("K6")
$c.Value = 6195
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(10).LineStyle = 1
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.ColorIndex = 1
$c.Interior.ThemeColor = 0
$c.VerticalAlignment = -4108
$c.NumberFormat = "#\ ##0"

# --- Widen the newly-used columns (K:O) to match the rest of the table (B:J) ---
$ws.Range("K1:O6").ColumnWidth = 7.83

"done"
